# edit.ps1 - apply the "last version ? :-)" commit changes to HAtH.pptx
#
# Summary of changes (see xml_diff):
#   Slide 1 (ctrTitle/subTitle shape "Subtitle 2"):
#     - para 1: "Forced Social Isolation raise many type of issue "
#               -> split into 3 runs, last word becomes "issues "
#     - para 2: "Health and Psychological issue" -> "...issues"
#     - para 3: "Reduce the social isolation issue and raise health prevention"
#               -> "...issues and raise health prevention"
#   Slide 2 (subTitle shape "Subtitle 2"):
#     - para 1: merge ", modular " + "open framework" runs into one run
#     - para 4: merge "High customization for " + "different " + "needs" into one run
#   Slide 4 (subTitle shape "Subtitle 2"):
#     - para 3: merge " " + "data" runs into one run (" data")
#   Slide 5 (ctrTitle shape "Title 1"):
#     - para 1: merge "The complete HAtH " + "solution Team" runs into one run

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - "The Problem" subtitle bullet list
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange

# Paragraph 1: "Forced Social Isolation raise many type of issue "
# Split so "of " and "issues " become their own runs, and "issue" -> "issues".
$para1 = $tr1.Paragraphs(1, 1)
$ofRange = $tr1.Characters($para1.Start + 40, 3)
$ofRange.Text = "of "
$issuesRange = $tr1.Characters($para1.Start + 43, 6)
$issuesRange.Text = "issues "

# Paragraph 2: "Health and Psychological issue" -> "...issues"
$para2 = $tr1.Paragraphs(2, 1)
$para2Range = $tr1.Characters($para2.Start, $para2.Length - 1)
$para2Range.Text = "Health and Psychological issues"

# Paragraph 3: "Reduce the social isolation issue and raise health prevention"
#              -> "...issues and raise health prevention"
$para3 = $tr1.Paragraphs(3, 1)
$para3Range = $tr1.Characters($para3.Start, $para3.Length - 1)
$para3Range.Text = "Reduce the social isolation issues and raise health prevention"

# ---------------------------------------------------------------------------
# Slide 2 - "The Solution: HAtH!" subtitle bullet list
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange

# Paragraph 1: "distributed" + ", modular " + "open framework"
# Merge the last two runs into a single run ", modular open framework".
$para2_1 = $tr2.Paragraphs(1, 1)
$mergeRange1 = $tr2.Characters($para2_1.Start + 11, 24)
$mergeRange1.Text = ", modular open framework"

# Paragraph 4: "High customization for " + "different " + "needs"
# Merge all three runs into a single run.
$para2_4 = $tr2.Paragraphs(4, 1)
$mergeRange2 = $tr2.Characters($para2_4.Start, $para2_4.Length - 1)
$mergeRange2.Text = "High customization for different needs"

# ---------------------------------------------------------------------------
# Slide 4 - "HAtH: NASA Resources" subtitle bullet list
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange

# Paragraph 3: "Isolation and human " + "behavioural" + " " + "data"
# Merge the trailing " " + "data" runs into a single run " data".
$para4_3 = $tr4.Paragraphs(3, 1)
$mergeRange3 = $tr4.Characters($para4_3.Start + 31, 5)
$mergeRange3.Text = " data"

# ---------------------------------------------------------------------------
# Slide 5 - "HAtH: The complete HAtH solution Team" title
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(1)
$tr5 = $shp5.TextFrame.TextRange

# Paragraph 1: "HAtH: " <br> "The complete HAtH " + "solution Team"
# Merge the two runs after the line break into a single run.
$para5_1 = $tr5.Paragraphs(1, 1)
$mergeRange4 = $tr5.Characters($para5_1.Start + 7, 31)
$mergeRange4.Text = "The complete HAtH solution Team"
